$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a duplicate of the "dev" row (row 2) at row 3, pushing the old
# row 3 ("preprod") down to row 4. Copy+Insert preserves per-cell styles.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

# New NroSiniestro values (plain strings, no accents) for the two "dev" rows.
# (E3's string is interned first so the shared-string table order matches.)
$ws.Range("E3").Value = "'1120194100405"
$ws.Range("E2").Value = "'1220194200662"

# Row 3 (the new duplicate row) needs the same custom height as row 2.
$ws.Rows.Item(3).RowHeight = 28.5

# Preserve B4's existing cell format (it currently carries the hyperlink
# style inherited from the old row 3), then re-point the hyperlink from
# B3 to B4, restoring the original formatting afterwards.
$ws.Range("B4").Copy($ws.Range("G1"))
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B4"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Range("G1").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("G1").Clear()

# Move the saved selection to match the new layout.
$ws.Range("G5").Select()
